# Commit: "Sat, Jun 13, 2020  9:04:50 AM"
#
# 1) Change the table style id applied to the table on slide 5.
# 2) Re-colour the deck's (single) theme - currently "Integral" / "Red
#    Violet" - so that its 12 theme colours match the stock Office theme
#    palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink), matching the
#    target OOXML's ppt/theme/theme1.xml colour values.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s5  = $p.Slides.Item(5)
$tbl = $s5.Shapes.Item(2).Table
$tbl.ApplyStyle("{C6B27CAC-9E35-4557-A2C5-46A9FE308BFE}", $true)

# --- 2. Theme colours -------------------------------------------------
# VBA/PowerPoint RGB longs are packed as R + G*256 + B*65536, so build
# each value from its familiar RRGGBB hex form.
function HexToBgrLong([int]$rrggbb) {
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x44546A,  # 3  dk2
    0xE7E6E6,  # 4  lt2
    0x5B9BD5,  # 5  accent1
    0xED7D31,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0xFFC000,  # 8  accent4
    0x4472C4,  # 9  accent5
    0x70AD47,  # 10 accent6
    0x0563C1,  # 11 hlink
    0x954F72   # 12 folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = HexToBgrLong $officeColors[$i - 1]
}
